# Generate Report for Handback
#
# This applies the "handback" update to the localization-status workbook:
#   - Status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" for every localized file row
#     (this automatically updates the Overview sheet too, since it shares
#     the same text).
#   - The "Latest Target File" (E) and "Latest Handback File" (F) columns
#     are now populated (they mirror the source markdown file and the
#     handoff xlf file respectively) with working hyperlinks, for every
#     data row on the zh-cn and de-de sheets.
#   - The "Latest Handback DateTime" (G) column is stamped with the new
#     handback timestamp for every data row on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# The Overview sheet mirrors the same "Status" text for each file (it
# shares the same string as the per-locale sheets), so update it too.
$overview = $wb.Worksheets.Item("Overview")
foreach ($row in @(2, 3)) {
    $overview.Cells.Item($row, 2).Value = $newStatus
    $overview.Cells.Item($row, 3).Value = $newStatus
}

$sheetHandbackTime = @{
    "zh-cn" = "2016-03-10 12:51:09"
    "de-de" = "2016-03-10 12:51:18"
}

# HyperLink-style blue underline color (matches the workbook's custom
# "HyperLink" cell style, RGB FF6495ED) expressed as a BGR integer the
# way the Excel object model expects for Font.Color / the RGB() macro.
$hyperlinkColor = 15570276

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Build a lookup of existing hyperlink addresses/display text keyed
    # by "row_col" so the new hyperlinks in columns E/F can reuse the
    # exact same targets as the corresponding A/C cells.
    $hlAddress = @{}
    $hlDisplay = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $r = $hl.Range
        $key = $r.Row.ToString() + "_" + $r.Column.ToString()
        $hlAddress[$key] = $hl.Address
        $hlDisplay[$key] = $hl.TextToDisplay
    }

    $newTime = $sheetHandbackTime[$sheetName]

    foreach ($row in @(2, 3)) {
        # Status column (B): "Ready for handoff" -> "Handed back: in sync with en-US"
        $ws.Cells.Item($row, 2).Value = $newStatus

        # Column A key "row_1" (Source File Name) -> column E (Latest Target File)
        $aKey = $row.ToString() + "_1"
        $aAddress = $hlAddress[$aKey]
        $aDisplay = $hlDisplay[$aKey]
        $eCell = $ws.Cells.Item($row, 5)
        $eCell.Value = $aDisplay
        $ws.Hyperlinks.Add($eCell, $aAddress, "", "", $aDisplay)
        $eCell.Font.Underline = $true
        $eCell.Font.Color = $hyperlinkColor

        # Column C key "row_3" (Latest Handoff File) -> column F (Latest Handback File)
        $cKey = $row.ToString() + "_3"
        $cAddress = $hlAddress[$cKey]
        $cDisplay = $hlDisplay[$cKey]
        $fCell = $ws.Cells.Item($row, 6)
        $fCell.Value = $cDisplay
        $ws.Hyperlinks.Add($fCell, $cAddress, "", "", $cDisplay)
        $fCell.Font.Underline = $true
        $fCell.Font.Color = $hyperlinkColor

        # Latest Handback DateTime column (G)
        $ws.Cells.Item($row, 7).Value = $newTime
    }
}
